$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1512539.4
$ws.Range("I92").Value = 651364.0600000001
$ws.Range("K92").Value = 651364.0600000001
$ws.Range("M92").Value = -650116.0600000001
$ws.Range("H100").Value = 6906.3335
$ws.Range("I100").Value = 3370.7144
$ws.Range("K100").Value = 3370.7144
$ws.Range("M100").Value = -2829.7144
$ws.Range("H115").Value = 1605.8334
$ws.Range("J115").Value = 1962.5
$ws.Range("L115").Value = 5887.5
$ws.Range("N115").Value = -9021.5
$ws.Range("H138").Value = 5387.086
$ws.Range("J138").Value = 3349.4783
$ws.Range("L138").Value = 10048.4349
$ws.Range("N138").Value = -20328.4349
$ws.Range("H141").Value = 16561.643
$ws.Range("I141").Value = 1688
$ws.Range("J141").Value = 20618.092
$ws.Range("K141").Value = 5064
$ws.Range("L141").Value = 61854.276
$ws.Range("M141").Value = 116
$ws.Range("N141").Value = -72214.276

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2620.8494
$ws.Range("I32").Value = 1173.9524
$ws.Range("K32").Value = 1173.9524
$ws.Range("M32").Value = -886.9523999999999
$ws.Range("H43").Value = 44323
$ws.Range("J43").Value = 44253.4
$ws.Range("L43").Value = 44253.4
$ws.Range("N43").Value = -44879.4
$ws.Range("H45").Value = 65282.688
$ws.Range("I45").Value = 65282.688
$ws.Range("K45").Value = 65282.688
$ws.Range("M45").Value = -64905.688
$ws.Range("H63").Value = 19692.129
$ws.Range("J63").Value = 23619.28
$ws.Range("L63").Value = 23619.28
$ws.Range("N63").Value = -24991.28
$ws.Range("H66").Value = 19692.129
$ws.Range("J66").Value = 23619.28
$ws.Range("L66").Value = 118096.4
$ws.Range("N66").Value = -124960.4
$ws.Range("H132").Value = 2892.8667
$ws.Range("I132").Value = 2378.9524
$ws.Range("K132").Value = 7136.8572
$ws.Range("M132").Value = -4606.8572
$ws.Range("H137").Value = 129885
$ws.Range("J137").Value = 129885
$ws.Range("L137").Value = 129885
$ws.Range("N137").Value = -140085

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3666.6667
$ws.Range("I22").Value = 3000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2827
$ws.Range("H99").Value = 9583.634
$ws.Range("I99").Value = 7932.6665
$ws.Range("J99").Value = 16187.5
$ws.Range("K99").Value = 7932.6665
$ws.Range("L99").Value = 16187.5
$ws.Range("M99").Value = -6434.6665
$ws.Range("N99").Value = -19183.5
$ws.Range("H107").Value = 9658.406999999999
$ws.Range("I107").Value = 11686.571
$ws.Range("J107").Value = 2559.8333
$ws.Range("K107").Value = 11686.571
$ws.Range("L107").Value = 2559.8333
$ws.Range("M107").Value = -9766.571
$ws.Range("N107").Value = -6399.8333
$ws.Range("H134").Value = 20932712
$ws.Range("I134").Value = 2166.2646
$ws.Range("J134").Value = 100003660
$ws.Range("K134").Value = 6498.793799999999
$ws.Range("L134").Value = 300010980
$ws.Range("M134").Value = -3963.793799999999
$ws.Range("N134").Value = -300016050

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 287
$ws.Range("I7").Value = 219.25
$ws.Range("J7").Value = 388.625
$ws.Range("K7").Value = 219.25
$ws.Range("L7").Value = 388.625
$ws.Range("M7").Value = -106.25
$ws.Range("N7").Value = -614.625
$ws.Range("H105").Value = 1230.5454
$ws.Range("I105").Value = 1295.9445
$ws.Range("K105").Value = 1295.9445
$ws.Range("M105").Value = 451.0554999999999
$ws.Range("H107").Value = 1155.9736
$ws.Range("I107").Value = 1133.4231
$ws.Range("J107").Value = 1204.8334
$ws.Range("K107").Value = 1133.4231
$ws.Range("L107").Value = 1204.8334
$ws.Range("M107").Value = 786.5769
$ws.Range("N107").Value = -5044.8334
$ws.Range("H132").Value = 11113996
$ws.Range("J132").Value = 25645174
$ws.Range("L132").Value = 76935522
$ws.Range("N132").Value = -76940582

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 996.26666
$ws.Range("I12").Value = 968.5
$ws.Range("J12").Value = 1000.53845
$ws.Range("K12").Value = 2905.5
$ws.Range("L12").Value = 3001.61535
$ws.Range("M12").Value = -2732.5
$ws.Range("N12").Value = -3347.61535
$ws.Range("H33").Value = 16750864
$ws.Range("I33").Value = 78
$ws.Range("J33").Value = 25126256
$ws.Range("K33").Value = 468
$ws.Range("L33").Value = 150757536
$ws.Range("M33").Value = -185
$ws.Range("N33").Value = -150758102
$ws.Range("H87").Value = 12699.3
$ws.Range("I87").Value = 7429.7144
$ws.Range("K87").Value = 22289.1432
$ws.Range("M87").Value = -21041.1432
$ws.Range("H90").Value = 12699.3
$ws.Range("I90").Value = 7429.7144
$ws.Range("K90").Value = 66867.4296
$ws.Range("M90").Value = -60627.4296
$ws.Range("J122").Value = 15873182
$ws.Range("L122").Value = 142858638
$ws.Range("N122").Value = -142863538
$ws.Range("H132").Value = 1271.0714
$ws.Range("I132").Value = 1058
$ws.Range("J132").Value = 1430.875
$ws.Range("K132").Value = 9522
$ws.Range("L132").Value = 12877.875
$ws.Range("M132").Value = -6992
$ws.Range("N132").Value = -17937.875
$ws.Range("H140").Value = 3215.739
$ws.Range("I140").Value = 1997.8235
$ws.Range("J140").Value = 6666.5
$ws.Range("K140").Value = 5993.470499999999
$ws.Range("L140").Value = 19999.5
$ws.Range("M140").Value = -813.4704999999994
$ws.Range("N140").Value = -30359.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 85.17646999999999
$ws.Range("I2").Value = 62.166668
$ws.Range("J2").Value = 140.4
$ws.Range("K2").Value = 62.166668
$ws.Range("L2").Value = 140.4
$ws.Range("M2").Value = 50.833332
$ws.Range("N2").Value = -366.4
$ws.Range("H44").Value = 44888.8
$ws.Range("J44").Value = 44888.8
$ws.Range("L44").Value = 44888.8
$ws.Range("N44").Value = -46080.8
$ws.Range("H62").Value = 46966.332
$ws.Range("I62").Value = 45449.5
$ws.Range("J62").Value = 50000
$ws.Range("K62").Value = 45449.5
$ws.Range("L62").Value = 50000
$ws.Range("M62").Value = -44763.5
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 46966.332
$ws.Range("I65").Value = 45449.5
$ws.Range("J65").Value = 50000
$ws.Range("K65").Value = 136348.5
$ws.Range("L65").Value = 150000
$ws.Range("M65").Value = -132916.5
$ws.Range("N65").Value = -156864
$ws.Range("H70").Value = 4882.778
$ws.Range("I70").Value = 4927.2856
$ws.Range("J70").Value = 4727
$ws.Range("K70").Value = 4927.2856
$ws.Range("L70").Value = 4727
$ws.Range("M70").Value = -4657.2856
$ws.Range("N70").Value = -5267
$ws.Range("H73").Value = 4882.778
$ws.Range("I73").Value = 4927.2856
$ws.Range("J73").Value = 4727
$ws.Range("K73").Value = 4927.2856
$ws.Range("L73").Value = 4727
$ws.Range("M73").Value = -3991.2856
$ws.Range("N73").Value = -6599
$ws.Range("H97").Value = 1515.6666
$ws.Range("I97").Value = 1025.2222
$ws.Range("K97").Value = 1025.2222
$ws.Range("M97").Value = -529.2221999999999
$ws.Range("H113").Value = 4669.6665
$ws.Range("I113").Value = 4505
$ws.Range("K113").Value = 4505
$ws.Range("M113").Value = -2335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 5000000
$ws.Range("I23").Value = 5000000
$ws.Range("K23").Value = 5000000
$ws.Range("M23").Value = -4999770
$ws.Range("H55").Value = 1139.6451
$ws.Range("I55").Value = 1208.6111
$ws.Range("J55").Value = 1044.1538
$ws.Range("K55").Value = 1208.6111
$ws.Range("L55").Value = 1044.1538
$ws.Range("M55").Value = -1035.6111
$ws.Range("N55").Value = -1390.1538
$ws.Range("H61").Value = 2516.8462
$ws.Range("I61").Value = 1816.579
$ws.Range("J61").Value = 4417.5713
$ws.Range("K61").Value = 1816.579
$ws.Range("L61").Value = 4417.5713
$ws.Range("M61").Value = -1614.579
$ws.Range("N61").Value = -4821.5713
$ws.Range("H82").Value = 2778.8
$ws.Range("I82").Value = 3150
$ws.Range("J82").Value = 2222
$ws.Range("K82").Value = 3150
$ws.Range("L82").Value = 2222
$ws.Range("M82").Value = -2789
$ws.Range("N82").Value = -2944
$ws.Range("H85").Value = 2778.8
$ws.Range("I85").Value = 3150
$ws.Range("J85").Value = 2222
$ws.Range("K85").Value = 3150
$ws.Range("L85").Value = 2222
$ws.Range("M85").Value = -1902
$ws.Range("N85").Value = -4718
$ws.Range("H113").Value = 2516.8462
$ws.Range("I113").Value = 1816.579
$ws.Range("J113").Value = 4417.5713
$ws.Range("K113").Value = 1816.579
$ws.Range("L113").Value = 4417.5713
$ws.Range("M113").Value = 353.421
$ws.Range("N113").Value = -8757.5713
$ws.Range("H122").Value = 3167.2727
$ws.Range("I122").Value = 2985
$ws.Range("K122").Value = 8955
$ws.Range("M122").Value = -6505

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 59999
$ws.Range("J54").Value = 59999
$ws.Range("L54").Value = 59999
$ws.Range("N54").Value = -61039
$ws.Range("H75").Value = 42142.57
$ws.Range("I75").Value = 29999.5
$ws.Range("J75").Value = 46999.8
$ws.Range("K75").Value = 29999.5
$ws.Range("L75").Value = 46999.8
$ws.Range("M75").Value = -29063.5
$ws.Range("N75").Value = -48871.8
$ws.Range("H78").Value = 42142.57
$ws.Range("I78").Value = 29999.5
$ws.Range("J78").Value = 46999.8
$ws.Range("K78").Value = 89998.5
$ws.Range("L78").Value = 140999.4
$ws.Range("M78").Value = -85318.5
$ws.Range("N78").Value = -150359.4
$ws.Range("H100").Value = 654859.6
$ws.Range("I100").Value = 727426.7
$ws.Range("K100").Value = 1454853.4
$ws.Range("M100").Value = -1454312.4
$ws.Range("H130").Value = 57499.5
$ws.Range("J130").Value = 57499.5
$ws.Range("L130").Value = 57499.5
$ws.Range("N130").Value = -67539.5
$ws.Range("H132").Value = 2112.1738
$ws.Range("I132").Value = 1408.9166
$ws.Range("J132").Value = 2879.3635
$ws.Range("K132").Value = 4226.7498
$ws.Range("L132").Value = 8638.0905
$ws.Range("M132").Value = -1696.7498
$ws.Range("N132").Value = -13698.0905
$ws.Range("H136").Value = 7001.1304
$ws.Range("I136").Value = 13790.143
$ws.Range("J136").Value = 4030.9375
$ws.Range("K136").Value = 41370.429
$ws.Range("L136").Value = 12092.8125
$ws.Range("M136").Value = -38820.429
$ws.Range("N136").Value = -17192.8125
